$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Insert a new row above row 6 (shifts existing rows 6..24 down to 7..25)
$ws.Rows.Item(6).Insert()

# Fill in the new row's data (Kidney / Solid Renal Mass / Clip 1 B-mode + Color / link text)
$ws.Range("A6").Value = "Kidney"
$ws.Range("B6").Value = "Solid Renal Mass "
$ws.Range("C6").Value = "Clip 1 B-mode + Color"
$ws.Range("D6").Value = "https://youtu.be/dGOyCO-pP4g"

# D6 uses the "hyperlink-looking" style (style index 1 in styles.xml) even
# though it is not an actual clickable hyperlink object.
$ws.Range("D6").Style = $ws.Range("D5").Style

# Update the active selection like Excel would after this edit
$ws.Range("A6").Select()
